# Weekly update: insert two new price records for "Pepino dulce" (Cultivar IV
# Región) at the top of the data table (rows 97-98), pushing the existing
# records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 97 (shifts 97:219 -> 99:221)
$ws.Rows("97:98").Insert()

# --- New row 97 ---------------------------------------------------------
$ws.Cells.Item(97, 1).Value = 10
$ws.Cells.Item(97, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(97, 3).Value = "La Araucanía"
$ws.Cells.Item(97, 4).Value = 44679
$ws.Cells.Item(97, 5).Value = 9
$ws.Cells.Item(97, 6).Value = 100112043
$ws.Cells.Item(97, 7).Value = "Pepino dulce"
$ws.Cells.Item(97, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(97, 9).Value = "Extra"
$ws.Cells.Item(97, 10).Value = 150
$ws.Cells.Item(97, 11).Value = 19000
$ws.Cells.Item(97, 12).Value = 20000
$ws.Cells.Item(97, 13).Value = 19333
$ws.Cells.Item(97, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(97, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(97, 16).Value = 1074
$ws.Cells.Item(97, 17).Value = 18
$ws.Cells.Item(97, 18).Value = "Hortaliza"

# --- New row 98 ---------------------------------------------------------
$ws.Cells.Item(98, 1).Value = 10
$ws.Cells.Item(98, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(98, 3).Value = "La Araucanía"
$ws.Cells.Item(98, 4).Value = 44679
$ws.Cells.Item(98, 5).Value = 9
$ws.Cells.Item(98, 6).Value = 100112043
$ws.Cells.Item(98, 7).Value = "Pepino dulce"
$ws.Cells.Item(98, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 200
$ws.Cells.Item(98, 11).Value = 17000
$ws.Cells.Item(98, 12).Value = 17000
$ws.Cells.Item(98, 13).Value = 17000
$ws.Cells.Item(98, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(98, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(98, 16).Value = 944
$ws.Cells.Item(98, 17).Value = 18
$ws.Cells.Item(98, 18).Value = "Hortaliza"
